$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text-formatted values (e.g. "1.00.01"-style
# thousands-dotted prices, or percentages with padding spaces) in the source data,
# so make sure Excel keeps them as text instead of auto-converting to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '26.988.09'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '1.563.65'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '207.74'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '22.12'
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').Value = '0.0600'
$ws.Range('E10').Value = '  +2.16%  '
$ws.Range('D11').Value = '0.0857'
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').Value = '1.785.99'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '1.561.16'
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').Value = '3.76'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = '0.520'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '62.03'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = '26.981.68'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('D19').Value = '216.10'
$ws.Range('E19').Value = '  -1.45%  '
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').Value = '4.11'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('D23').Value = '9.21'
$ws.Range('E23').Value = '  -0.66%  '
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').Value = '152.99'
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('E31').Value = '  +1.27%  '
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('D34').Value = '1.422.90'
$ws.Range('E34').Value = '  -1.62%  '
$ws.Range('E35').Value = '  +2.87%  '
$ws.Range('D36').Value = '1.07'
$ws.Range('E36').Value = '  +10.90%  '
$ws.Range('E37').Value = '  +1.89%  '
$ws.Range('E38').Value = '  -0.40%  '
$ws.Range('D39').Value = '0.533'
$ws.Range('E39').Value = '  +1.82%  '
$ws.Range('D40').Value = '5.81'
$ws.Range('E40').Value = '  +1.55%  '
$ws.Range('D41').Value = '0.808'
$ws.Range('E41').Value = '  -0.81%  '
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '2.32'
$ws.Range('E43').Value = '  +1.98%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '1.01'
$ws.Range('E44').Value = '  +1.86%  '
$ws.Range('D45').Value = '64.75'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('E46').Value = '  -1.64%  '
$ws.Range('D47').Value = '1.698.62'
$ws.Range('E47').Value = '  -0.23%  '
$ws.Range('D48').Value = '87.33'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('D49').Value = '0.0520'
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0101'
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.0960'
$ws.Range('E51').Value = '  -0.49%  '
